# Apply the changes described by the diff:
#   1. Drop the forced "fixed" table layout on both tables in the body
#      (i.e. let Word auto-fit the table again instead of pinning it to
#      fixed column widths) -- removes <w:tblLayout w:type="fixed"/>.
#   2. Give the "Abstract" style the same 300-twip (15pt) space-before
#      that "Abstract Title" used to have, now that the title style is
#      about to go away (w:before 100 -> 300 in <w:spacing>).
#   3. Remove the unused custom paragraph style "Abstract Title"
#      (styleId "AbstractTitle").
#
# NOTE: style/table handles are positional and get reseated by a
# mutating call (e.g. Delete), so do the rename-sensitive edits first
# and look styles back up by name only while the collection is still
# untouched.

$d = $word.ActiveDocument

# 1) Tables: allow auto-fit again on every table in the document instead
#    of forcing a fixed layout.
$tableCount = $d.Tables.Count
foreach ($t in $d.Tables) {
    $t.AllowAutoFit = $true
}

# 2) Bump "Abstract" style's space-before from 5pt (100 twips) to
#    15pt (300 twips) to match the "after" spacing (300/300). Do this
#    before deleting "Abstract Title" below.
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 15
$abstractSpaceBefore = $abstract.ParagraphFormat.SpaceBefore

# 3) Delete the now-unused "Abstract Title" paragraph style.
$abstractTitle = $d.Styles("AbstractTitle")
$deletedStyleName = $abstractTitle.NameLocal
$abstractTitle.Delete()

Write-Output "Tables set to auto-fit: $tableCount"
Write-Output "Abstract SpaceBefore now: $abstractSpaceBefore"
Write-Output "Deleted style: $deletedStyleName"
